$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace "Chikkaballapura" with "Chikballapur" in column G for data rows 3-69,
# except rows 7, 19 and 30 (which keep their original District text as-is:
# row 7/19 have a different district value already, row 30 has different casing).
for ($r = 3; $r -le 69; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $v = $cell.Value2
    if ($v -ne $null -and $v.Equals("Chikkaballapura")) {
        $cell.Value2 = "Chikballapur"
    }
}

# Clear the empty inline-string placeholders in F7 and F19 so the cells no
# longer exist in the saved XML.
$ws.Cells.Item(7, 6).Clear()
$ws.Cells.Item(19, 6).Clear()
